$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 42
$ws1.Range("F4").Value = 38
$ws1.Range("F6").Value = 5303
$ws1.Range("F7").Value = 180
$ws1.Range("F8").Value = 104
$ws1.Range("F10").Value = 363
$ws1.Range("F12").Value = 68

# Sheet "全部类型" (all types) - same underlying events, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 42
$ws4.Range("F7").Value = 38
$ws4.Range("F9").Value = 5303
$ws4.Range("F10").Value = 180
$ws4.Range("F11").Value = 104
$ws4.Range("F14").Value = 363
$ws4.Range("F16").Value = 68
